$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-13) holds a "Förändrad" date that was bumped by one day
# (serial 46061 -> 46062, i.e. 2026-02-08 -> 2026-02-09) for every data row.
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
